$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D:E data range to Text format so that numeric-looking
# strings (prices, percentages) are NOT auto-converted to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '41.957.58'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.229.58'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '251.02'
$ws.Range('E5').Value = '  +7.93%  '
$ws.Range('D6').Value = '0.632'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '71.21'
$ws.Range('E7').Value = '  +3.22%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = '0.594'
$ws.Range('E9').Value = '  +6.38%  '
$ws.Range('D10').Value = '41.28'
$ws.Range('E10').Value = '  +15.43%  '
$ws.Range('D11').Value = '0.0971'
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('D12').Value = '58.33'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('D13').Value = '7.24'
$ws.Range('E13').Value = '  +7.81%  '
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').Value = '2.558.70'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').Value = '15.01'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '0.867'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').Value = '2.229.61'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').Value = '41.819.29'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').Value = '0.0₃0976'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = '6.23'
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').Value = '73.04'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = '235.75'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  +8.30%  '
$ws.Range('D25').Value = '4.22'
$ws.Range('E25').Value = '  +15.40%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '2.56'
$ws.Range('E27').Value = '  +8.76%  '
$ws.Range('D28').Value = '10.78'
$ws.Range('E28').Value = '  +7.83%  '
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').Value = '171.11'
$ws.Range('E30').Value = '  +1.44%  '
$ws.Range('D31').Value = '20.81'
$ws.Range('E31').Value = '  +1.18%  '
$ws.Range('E32').Value = '  +2.97%  '
$ws.Range('D33').Value = '0.126'
$ws.Range('E33').Value = '  -0.71%  '
$ws.Range('D34').Value = '5.62'
$ws.Range('E34').Value = '  +7.08%  '
$ws.Range('D35').Value = '0.0729'
$ws.Range('E35').Value = '  +2.51%  '
$ws.Range('D36').Value = '4.72'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').Value = '26.05'
$ws.Range('E37').Value = '  +19.37%  '
$ws.Range('D38').Value = '3.93'
$ws.Range('E38').Value = '  +8.91%  '
$ws.Range('D39').Value = '0.0301'
$ws.Range('E39').Value = '  +13.98%  '
$ws.Range('D40').Value = '2.31'
$ws.Range('E40').Value = '  +2.80%  '
$ws.Range('D41').Value = '5.98'
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').Value = '67.92'
$ws.Range('E42').Value = '  +2.60%  '
$ws.Range('E43').Value = '  +10.52%  '
$ws.Range('D44').Value = '11.83'
$ws.Range('E44').Value = '  +18.28%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').Value = '4.91'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '8.83'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').Value = '0.103'
$ws.Range('E47').Value = '  +2.50%  '
$ws.Range('D48').Value = '4.71'
$ws.Range('E48').Value = '  +8.62%  '
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('E50').Value = '  +7.77%  '
$ws.Range('B51').Value = 'BitTorrent-New'
$ws.Range('C51').Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range('D51').Value = '0.0₃0149'
$ws.Range('E51').Value = '  +6.08%  '

# Restore original (default/general) style so no stray number format
# is left behind on the cells.
$dataRange.Style = "Normal"
